$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (prices in column D, 1h volume % in column E).
# A leading apostrophe forces Excel to store numeric-looking values
# (e.g. "1.000", "0.7780") as text, matching the inlineStr cells in the
# original workbook instead of letting Excel coerce them to numbers.
$ws.Range("D2").Value = "26.428.91"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "1.847.31"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'264.96"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.5214"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").Value = "'0.3276"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "'0.06814"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "'18.86"
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("D11").Value = "'0.7780"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "'0.07777"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "1.839.82"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "'88.18"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "'5.017"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'13.95"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "'0.000007990"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "26.459.15"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "2.070.35"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").Value = "'4.645"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'9.565"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").Value = "'6.006"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "'144.49"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("E26").Value = "  -7.76%  "
$ws.Range("D27").Value = "'1.671"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "'17.04"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "'112.22"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'4.182"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "'4.142"
$ws.Range("D32").Value = "'0.08764"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "'1.137"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "'0.7200"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'2.856"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "'3.105"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").Value = "'2.212"
$ws.Range("E39").Value = "  -3.50%  "
$ws.Range("D40").Value = "'0.4870"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").Value = "'0.9086"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "'111.37"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").Value = "'6.072"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "'0.9994"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "'7.717"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.05943"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4172"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").Value = "'9.087"
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("E49").Value = "  -5.96%  "
$ws.Range("D50").Value = "'34.97"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "'0.8925"
$ws.Range("E51").Value = "  +2.25%  "
